$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28: The Writing Is Not on the Wall / Leve Item ID 27772
$ws.Range("H28").Value = 1434.5555
$ws.Range("I28").Value = 1275.2727
$ws.Range("K28").Value = 1275.2727
$ws.Range("M28").Value = -790.2727

# Row 33: Glazed and Confused / Leve Item ID 5512
$ws.Range("H33").Value = 11111302
$ws.Range("I33").Value = 11111302
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 11111302
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -11111073
$ws.Range("N33").Value = $null

# Row 62: The Mustache Suits Him / Leve Item ID 27781
$ws.Range("H62").Value = 31906.553
$ws.Range("I62").Value = 44443.875
$ws.Range("K62").Value = 44443.875
$ws.Range("M62").Value = -43819.875

# Row 65: Forgery of Convenience (L) / Leve Item ID 27781
$ws.Range("H65").Value = 31906.553
$ws.Range("I65").Value = 44443.875
$ws.Range("K65").Value = 222219.375
$ws.Range("M65").Value = -219099.375

# Row 76: Warding Off Temptation / Leve Item ID 12602
$ws.Range("H76").Value = 4748.5
$ws.Range("I76").Value = 4247.5
$ws.Range("K76").Value = 4247.5
$ws.Range("M76").Value = -3932.5

# Row 79: The Garden of Arcane Delights (L) / Leve Item ID 12602
$ws.Range("H79").Value = 4748.5
$ws.Range("I79").Value = 4247.5
$ws.Range("K79").Value = 4247.5
$ws.Range("M79").Value = -3155.5

# Row 80: Cleansing the Wicked Humours / Leve Item ID 12605
$ws.Range("H80").Value = 2262.7896
$ws.Range("I80").Value = 2032.4445
$ws.Range("K80").Value = 6097.333500000001
$ws.Range("M80").Value = -5099.333500000001

# Row 83: Washing Away the Sins (L) / Leve Item ID 12605
$ws.Range("H83").Value = 2262.7896
$ws.Range("I83").Value = 2032.4445
$ws.Range("K83").Value = 18292.0005
$ws.Range("M83").Value = -13300.0005

# Row 112: Making Ends Meet / Leve Item ID 27960
$ws.Range("H112").Value = 92488.17999999999
$ws.Range("J112").Value = 126760.125
$ws.Range("L112").Value = 380280.375
$ws.Range("N112").Value = -382496.375

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff / Leve Item ID 43999
$ws.Range("H61").Value = 3664.0715
$ws.Range("I61").Value = 2913.182
$ws.Range("K61").Value = 2913.182
$ws.Range("M61").Value = -2701.182

# Row 88: The Mast Chance / Leve Item ID 12530
$ws.Range("H88").Value = 2916
$ws.Range("I88").Value = 3482.4
$ws.Range("J88").Value = 1500
$ws.Range("K88").Value = 3482.4
$ws.Range("L88").Value = 1500
$ws.Range("M88").Value = -3076.4
$ws.Range("N88").Value = -2312

# Row 91: The Rose and the Riveter (L) / Leve Item ID 12530
$ws.Range("H91").Value = 2916
$ws.Range("I91").Value = 3482.4
$ws.Range("J91").Value = 1500
$ws.Range("K91").Value = 3482.4
$ws.Range("L91").Value = 1500
$ws.Range("M91").Value = -2078.4
$ws.Range("N91").Value = -4308

# Row 122: Haste for High Durium / Leve Item ID 36168
$ws.Range("H122").Value = 1940.6154
$ws.Range("I122").Value = 1602.125
$ws.Range("J122").Value = 2482.2
$ws.Range("K122").Value = 4806.375
$ws.Range("L122").Value = 7446.599999999999
$ws.Range("M122").Value = -2356.375
$ws.Range("N122").Value = -12346.6

# Row 136: Metal with Mettle / Leve Item ID 43999
$ws.Range("H136").Value = 3664.0715
$ws.Range("I136").Value = 2913.182
$ws.Range("K136").Value = 8739.545999999998
$ws.Range("M136").Value = -6189.545999999998

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin / Leve Item ID 12526
$ws.Range("H86").Value = 2030.4736
$ws.Range("I86").Value = 1140.9
$ws.Range("J86").Value = 3018.889
$ws.Range("K86").Value = 1140.9
$ws.Range("L86").Value = 3018.889
$ws.Range("M86").Value = -17.90000000000009
$ws.Range("N86").Value = -5264.889

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Leve Item ID 12526
$ws.Range("H89").Value = 2030.4736
$ws.Range("I89").Value = 1140.9
$ws.Range("J89").Value = 3018.889
$ws.Range("K89").Value = 5704.5
$ws.Range("L89").Value = 15094.445
$ws.Range("M89").Value = -88.5
$ws.Range("N89").Value = -26326.445

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof / Leve Item ID 27691
$ws.Range("H16").Value = 1471.5
$ws.Range("I16").Value = 1258.5714
$ws.Range("K16").Value = 1258.5714
$ws.Range("M16").Value = -971.5714

# Row 31: Wall Not Found / Leve Item ID 44023
$ws.Range("H31").Value = 29531.71
$ws.Range("I31").Value = 38477.63
$ws.Range("J31").Value = 7573.5454
$ws.Range("K31").Value = 38477.63
$ws.Range("L31").Value = 7573.5454
$ws.Range("M31").Value = -38182.63
$ws.Range("N31").Value = -8163.5454

# Row 34: Armoires of the Rich and Famous / Leve Item ID 44023
$ws.Range("H34").Value = 29531.71
$ws.Range("I34").Value = 38477.63
$ws.Range("J34").Value = 7573.5454
$ws.Range("K34").Value = 38477.63
$ws.Range("L34").Value = 7573.5454
$ws.Range("M34").Value = -38275.63
$ws.Range("N34").Value = -7977.5454

# Row 99: O Pine / Leve Item ID 36198
$ws.Range("H99").Value = 3271.9443
$ws.Range("I99").Value = 2926.4
$ws.Range("J99").Value = 4999.6665
$ws.Range("K99").Value = 2926.4
$ws.Range("L99").Value = 4999.6665
$ws.Range("M99").Value = -1428.4
$ws.Range("N99").Value = -7995.6665

# Row 113: Patient Patients / Leve Item ID 27691
$ws.Range("H113").Value = 1471.5
$ws.Range("I113").Value = 1258.5714
$ws.Range("K113").Value = 1258.5714
$ws.Range("M113").Value = 911.4286

# Row 126: A Better Conductor / Leve Item ID 36198
$ws.Range("H126").Value = 3271.9443
$ws.Range("I126").Value = 2926.4
$ws.Range("J126").Value = 4999.6665
$ws.Range("K126").Value = 8779.200000000001
$ws.Range("L126").Value = 14998.9995
$ws.Range("M126").Value = -6309.200000000001
$ws.Range("N126").Value = -19938.9995

$ws = $wb.Worksheets.Item("GSM")
# Row 18: Gorgeous Gorget / Leve Item ID 4309
$ws.Range("H18").Value = 55558100
$ws.Range("I18").Value = 55558100
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 55558100
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -55557807
$ws.Range("N18").Value = $null

# Row 25: I Am a Rock / Leve Item ID 2687
$ws.Range("H25").Value = 2890.5715
$ws.Range("I25").Value = 1693.75
$ws.Range("J25").Value = 4486.3335
$ws.Range("K25").Value = 1693.75
$ws.Range("L25").Value = 4486.3335
$ws.Range("M25").Value = -1164.75
$ws.Range("N25").Value = -5544.3335

# Row 44: Actually, It''s Loyalty / Leve Item ID 4143
$ws.Range("H44").Value = 29500.5
$ws.Range("J44").Value = 29500.5
$ws.Range("L44").Value = 29500.5
$ws.Range("N44").Value = -30692.5

# Row 48: Dead Can''t Defang / Leve Item ID 4337
$ws.Range("H48").Value = 15000
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").Value = $null

# Row 122: Awarding Academic Excellence / Leve Item ID 36182
$ws.Range("H122").Value = 840
$ws.Range("I122").Value = 840
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2520
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -70
$ws.Range("N122").Value = $null

# Row 126: Gold Rush Order / Leve Item ID 36184
$ws.Range("H126").Value = 25638.857
$ws.Range("I126").Value = 37716
$ws.Range("J126").Value = 3900
$ws.Range("K126").Value = 113148
$ws.Range("L126").Value = 11700
$ws.Range("M126").Value = -110678
$ws.Range("N126").Value = -16640

$ws = $wb.Worksheets.Item("LTW")
# Row 46: Supply Side Logic / Leve Item ID 5282
$ws.Range("H46").Value = 1850
$ws.Range("I46").Value = 1633.3334
$ws.Range("J46").Value = 2500
$ws.Range("K46").Value = 1633.3334
$ws.Range("L46").Value = 2500
$ws.Range("M46").Value = -1445.3334
$ws.Range("N46").Value = -2876

# Row 93: Hide to Go Seek / Leve Item ID 19993
$ws.Range("H93").Value = 3297.4707
$ws.Range("I93").Value = 3119.8462
$ws.Range("K93").Value = 3119.8462
$ws.Range("M93").Value = -1871.8462

$ws = $wb.Worksheets.Item("WVR")
# Row 129: Lifetime of Gleaning / Leve Item ID 35429
$ws.Range("H129").Value = 34475
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 34475
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 34475
$ws.Range("M129").Value = $null
$ws.Range("N129").Value = -44475

